# The "Import Files" (C) and "Export Files" (D) columns were off by one row
# (rows 11-20) relative to their correct Job Name/Schedule/Status/Comment
# rows. This shifts the C/D content down by one row (row 20 picks up the
# values already present in row 22, the "DB server Disk Space Monitoring"
# summary row, which itself stays untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 11
$endRow = 20

# Capture the original values first so the shift doesn't clobber data we
# still need to read further down the chain.
$cValues = @{}
$dValues = @{}
for ($r = $startRow; $r -le ($endRow + 2); $r++) {
    $cValues[$r] = $ws.Cells.Item($r, 3).Value2
    $dValues[$r] = $ws.Cells.Item($r, 4).Value2
}

for ($r = $startRow; $r -le $endRow; $r++) {
    if ($r -eq $endRow) {
        # Row 20 <- Row 22 (row 21 is a blank separator row)
        $ws.Cells.Item($r, 3).Value2 = $cValues[$r + 2]
        $ws.Cells.Item($r, 4).Value2 = $dValues[$r + 2]
    } else {
        # Row N <- Row N+1
        $ws.Cells.Item($r, 3).Value2 = $cValues[$r + 1]
        $ws.Cells.Item($r, 4).Value2 = $dValues[$r + 1]
    }
}
